$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "29.821.06"
$ws.Range("E2").Value = "  +4.88%  "

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.612.15"
$ws.Range("E3").Value = "  +3.80%  "

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.995"
$ws.Range("E4").Value = "  -0.44%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "213.62"
$ws.Range("E5").Value = "  +1.40%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.516"
$ws.Range("E6").Value = "  +6.84%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.994"
$ws.Range("E7").Value = "  -0.58%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "26.85"
$ws.Range("E8").Value = "  +11.37%  "

# Row 9
$ws.Range("E9").Value = "  +3.07%  "

# Row 10
$ws.Range("E10").Value = "  +2.46%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0911"
$ws.Range("E11").Value = "  +2.44%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.842.01"
$ws.Range("E12").Value = "  +3.73%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.610.47"
$ws.Range("E13").Value = "  +3.65%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "29.855.31"
$ws.Range("E14").Value = "  +4.93%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.537"
$ws.Range("E15").Value = "  +5.27%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.76"
$ws.Range("E16").Value = "  +3.76%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "244.70"
$ws.Range("E17").Value = "  +6.85%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "63.50"

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.0₃0695"
$ws.Range("E20").Value = "  +3.33%  "

# Row 21
$ws.Range("E21").Value = "  -0.49%  "

# Row 22
$ws.Range("E22").Value = "  +4.34%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "9.26"
$ws.Range("E23").Value = "  +3.99%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.09"
$ws.Range("E24").Value = "  +3.80%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "156.36"
$ws.Range("E25").Value = "  +3.66%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "15.34"
$ws.Range("E26").Value = "  +4.11%  "

# Row 27
$ws.Range("E27").Value = "  +5.41%  "

# Row 28
$ws.Range("E28").Value = "  +2.91%  "

# Row 29
$ws.Range("E29").Value = "  -0.46%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.0474"
$ws.Range("E30").Value = "  +1.40%  "

# Row 31
$ws.Range("E31").Value = "  +0.99%  "

# Row 32
$ws.Range("E32").Value = "  +2.88%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.442.43"
$ws.Range("E33").Value = "  +4.23%  "

# Row 34
$ws.Range("E34").Value = "  +3.71%  "

# Row 35
$ws.Range("E35").Value = "  -0.25%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.51"

# Row 37
$ws.Range("E37").Value = "  +9.81%  "

# Row 38
$ws.Range("E38").Value = "  +0.54%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.0167"
$ws.Range("E39").Value = "  +3.30%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.537"
$ws.Range("E40").Value = "  +5.07%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "55.22"
$ws.Range("E41").Value = "  +27.84%  "

# Row 42
$ws.Range("E42").Value = "  +1.40%  "

# Row 43
$ws.Range("E43").Value = "  +3.36%  "

# Row 44
$ws.Range("E44").Value = "  -0.50%  "

# Row 45
$ws.Range("E45").Value = "  +1.54%  "

# Row 46
$ws.Range("E46").Value = "  +7.02%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "5.30"
$ws.Range("E47").Value = "  -1.00%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.753.56"
$ws.Range("E48").Value = "  +4.03%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "86.98"
$ws.Range("E49").Value = "  +2.35%  "

# Row 50
$ws.Range("E50").Value = "  -4.26%  "

# Row 51
$ws.Range("B51").Value = "BabyDogeCoin"
$ws.Range("C51").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0₆0104"
$ws.Range("E51").Value = "  +2.27%  "
